# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Binance" conversion rate lines ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = [string][char]0x2705 + " 1000 Bs = 5.45 = 21701.32 pesos"
$newLine1 = [string][char]0x2705 + " 1000 Bs = 5.38 = 21377.75 pesos"
$oldLine2 = [string][char]0x2705 + " 21701.32 pesos = 5.43 = 960.89 Bs"
$newLine2 = [string][char]0x2705 + " 21377.75 pesos = 5.33 = 960.16 Bs"

$text = $ws1.Range("A1").Value
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$ws1.Range("A1").Value = $text

# --- tasas!N10/O10/N12/O12: update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 185.941
$ws2.Range("O10").Value = 3975
$ws2.Range("N12").Value = 4007.66
$ws2.Range("O12").Value = 180
